$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.164.09'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '1.604.59'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.16'
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  +1.03%  '
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.07'
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('D12').Value = '1.828.67'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '1.602.76'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('E14').Value = '  -1.17%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.508'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '26.153.24'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '60.59'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '198.98'
$ws.Range('E20').Value = '  +4.77%  '
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.40'
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('E23').Value = '  +0.87%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.131'
$ws.Range('E24').Value = '  +2.68%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '141.96'
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.75'
$ws.Range('E26').Value = '  +2.32%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.14'
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('E31').Value = '  +1.11%  '
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('E35').Value = '  -1.60%  '
$ws.Range('D36').Value = '1.107.64'
$ws.Range('E36').Value = '  +0.96%  '
$ws.Range('B37').Value = 'PaxDollar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.35'
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.501'
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.784'
$ws.Range('E41').Value = '  -0.93%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.782'
$ws.Range('E42').Value = '  +5.40%  '
$ws.Range('D43').Value = '1.742.75'
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.13'
$ws.Range('E44').Value = '  +1.29%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '92.95'
$ws.Range('E45').Value = '  -2.84%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.55'
$ws.Range('E46').Value = '  +7.76%  '
$ws.Range('D47').Value = '0.0₆0104'
$ws.Range('E47').Value = '  -7.31%  '
$ws.Range('E48').Value = '  +0.71%  '
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.409'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('E51').Value = '  -0.06%  '
